$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOCK_DATA")

$ws.Range("I1").Value = 'url'
$ws.Range("I1").Font.Bold = $true
$ws.Range("I2").Value = 'https://latimes.com/donec/vitae/nisi/nam/ultrices.json?convallis=sapien&morbi=iaculis&odio=congue&odio=vivamus&elementum=metus&eu=arcu&interdum=adipiscing&eu=molestie&tincidunt=hendrerit&in=at&leo=vulputate&maecenas=vitae&pulvinar=nisl&lobortis=aenean&est=lectus&phasellus=pellentesque&sit=eget&amet=nunc&erat=donec&nulla=quis&tempus=orci'
$ws.Range("I3").Value = 'http://omniture.com/sed/vestibulum/sit.html?nisl=pulvinar&aenean=nulla&lectus=pede&pellentesque=ullamcorper&eget=augue&nunc=a&donec=suscipit&quis=nulla&orci=elit&eget=ac&orci=nulla&vehicula=sed&condimentum=vel&curabitur=enim&in=sit&libero=amet&ut=nunc&massa=viverra&volutpat=dapibus&convallis=nulla&morbi=suscipit&odio=ligula&odio=in&elementum=lacus&eu=curabitur&interdum=at&eu=ipsum&tincidunt=ac&in=tellus&leo=semper&maecenas=interdum&pulvinar=mauris&lobortis=ullamcorper'
$ws.Range("I4").Value = 'https://cam.ac.uk/purus/phasellus/in/felis.js?eget=quisque&vulputate=erat&ut=eros&ultrices=viverra&vel=eget&augue=congue&vestibulum=eget&ante=semper&ipsum=rutrum&primis=nulla&in=nunc&faucibus=purus&orci=phasellus&luctus=in&et=felis&ultrices=donec&posuere=semper&cubilia=sapien&curae=a&donec=libero&pharetra=nam&magna=dui&vestibulum=proin&aliquet=leo&ultrices=odio&erat=porttitor&tortor=id&sollicitudin=consequat&mi=in&sit=consequat&amet=ut&lobortis=nulla&sapien=sed&sapien=accumsan&non=felis&mi=ut&integer=at&ac=dolor&neque=quis&duis=odio&bibendum=consequat&morbi=varius&non=integer&quam=ac&nec=leo&dui=pellentesque&luctus=ultrices&rutrum=mattis&nulla=odio&tellus=donec&in=vitae&sagittis=nisi&dui=nam&vel=ultrices&nisl=libero&duis=non&ac=mattis&nibh=pulvinar&fusce=nulla&lacus=pede&purus=ullamcorper&aliquet=augue&at=a&feugiat=suscipit&non=nulla&pretium=elit&quis=ac&lectus=nulla&suspendisse=sed&potenti=vel&in=enim&eleifend=sit&quam=amet&a=nunc&odio=viverra&in=dapibus&hac=nulla&habitasse=suscipit&platea=ligula&dictumst=in&maecenas=lacus'
$ws.Range("I5").Value = 'https://163.com/blandit/nam/nulla/integer/pede/justo.aspx?nullam=pulvinar&varius=lobortis&nulla=est&facilisi=phasellus&cras=sit&non=amet&velit=erat&nec=nulla&nisi=tempus&vulputate=vivamus&nonummy=in&maecenas=felis&tincidunt=eu&lacus=sapien&at=cursus&velit=vestibulum&vivamus=proin&vel=eu&nulla=mi&eget=nulla&eros=ac&elementum=enim&pellentesque=in&quisque=tempor&porta=turpis&volutpat=nec&erat=euismod&quisque=scelerisque&erat=quam&eros=turpis&viverra=adipiscing&eget=lorem&congue=vitae&eget=mattis&semper=nibh&rutrum=ligula&nulla=nec&nunc=sem&purus=duis&phasellus=aliquam&in=convallis&felis=nunc&donec=proin&semper=at&sapien=turpis&a=a&libero=pede&nam=posuere&dui=nonummy&proin=integer&leo=non&odio=velit&porttitor=donec&id=diam&consequat=neque&in=vestibulum&consequat=eget&ut=vulputate&nulla=ut&sed=ultrices&accumsan=vel&felis=augue&ut=vestibulum&at=ante&dolor=ipsum&quis=primis&odio=in&consequat=faucibus&varius=orci&integer=luctus&ac=et&leo=ultrices&pellentesque=posuere&ultrices=cubilia&mattis=curae&odio=donec&donec=pharetra&vitae=magna&nisi=vestibulum&nam=aliquet&ultrices=ultrices&libero=erat&non=tortor&mattis=sollicitudin&pulvinar=mi&nulla=sit'
$ws.Range("I6").Value = 'http://wiley.com/in/faucibus.html?congue=nulla&vivamus=facilisi&metus=cras&arcu=non&adipiscing=velit&molestie=nec&hendrerit=nisi&at=vulputate&vulputate=nonummy&vitae=maecenas&nisl=tincidunt&aenean=lacus&lectus=at&pellentesque=velit&eget=vivamus&nunc=vel&donec=nulla&quis=eget&orci=eros&eget=elementum&orci=pellentesque&vehicula=quisque&condimentum=porta&curabitur=volutpat&in=erat&libero=quisque&ut=erat&massa=eros&volutpat=viverra&convallis=eget'
$ws.Range("I7").Value = 'http://booking.com/aenean/lectus/pellentesque/eget/nunc/donec.xml?mauris=fermentum&enim=justo&leo=nec&rhoncus=condimentum&sed=neque&vestibulum=sapien&sit=placerat&amet=ante&cursus=nulla&id=justo&turpis=aliquam&integer=quis&aliquet=turpis&massa=eget&id=elit&lobortis=sodales&convallis=scelerisque&tortor=mauris&risus=sit&dapibus=amet&augue=eros&vel=suspendisse&accumsan=accumsan&tellus=tortor&nisi=quis&eu=turpis&orci=sed&mauris=ante&lacinia=vivamus&sapien=tortor&quis=duis&libero=mattis&nullam=egestas&sit=metus&amet=aenean&turpis=fermentum&elementum=donec&ligula=ut&vehicula=mauris&consequat=eget&morbi=massa&a=tempor&ipsum=convallis&integer=nulla&a=neque&nibh=libero'
$ws.Range("I8").Value = 'http://vistaprint.com/est/congue/elementum/in/hac.js?lacus=aliquet&curabitur=massa&at=id&ipsum=lobortis&ac=convallis'
$ws.Range("I9").Value = 'http://shutterfly.com/tempus/sit/amet.aspx?ac=non&neque=mauris&duis=morbi&bibendum=non&morbi=lectus&non=aliquam&quam=sit&nec=amet&dui=diam&luctus=in&rutrum=magna&nulla=bibendum&tellus=imperdiet&in=nullam&sagittis=orci'
$ws.Range("I10").Value = 'https://creativecommons.org/molestie/nibh/in/lectus/pellentesque.png?praesent=in&lectus=hac&vestibulum=habitasse&quam=platea&sapien=dictumst&varius=morbi&ut=vestibulum&blandit=velit&non=id&interdum=pretium&in=iaculis&ante=diam&vestibulum=erat&ante=fermentum&ipsum=justo&primis=nec&in=condimentum&faucibus=neque&orci=sapien&luctus=placerat&et=ante&ultrices=nulla&posuere=justo&cubilia=aliquam&curae=quis&duis=turpis&faucibus=eget&accumsan=elit&odio=sodales&curabitur=scelerisque&convallis=mauris&duis=sit&consequat=amet&dui=eros&nec=suspendisse&nisi=accumsan&volutpat=tortor&eleifend=quis&donec=turpis&ut=sed&dolor=ante&morbi=vivamus&vel=tortor&lectus=duis&in=mattis&quam=egestas&fringilla=metus&rhoncus=aenean&mauris=fermentum&enim=donec&leo=ut&rhoncus=mauris'
$ws.Range("I11").Value = 'https://paypal.com/ligula/vehicula/consequat/morbi/a/ipsum.jsp?ipsum=justo&integer=pellentesque&a=viverra&nibh=pede&in=ac&quis=diam&justo=cras&maecenas=pellentesque&rhoncus=volutpat&aliquam=dui&lacus=maecenas&morbi=tristique&quis=est&tortor=et&id=tempus&nulla=semper&ultrices=est&aliquet=quam&maecenas=pharetra&leo=magna&odio=ac&condimentum=consequat&id=metus&luctus=sapien&nec=ut&molestie=nunc&sed=vestibulum&justo=ante&pellentesque=ipsum&viverra=primis&pede=in&ac=faucibus&diam=orci&cras=luctus&pellentesque=et&volutpat=ultrices&dui=posuere&maecenas=cubilia&tristique=curae&est=mauris&et=viverra&tempus=diam&semper=vitae&est=quam&quam=suspendisse&pharetra=potenti&magna=nullam&ac=porttitor&consequat=lacus&metus=at&sapien=turpis&ut=donec&nunc=posuere&vestibulum=metus&ante=vitae&ipsum=ipsum&primis=aliquam&in=non&faucibus=mauris&orci=morbi&luctus=non&et=lectus&ultrices=aliquam&posuere=sit&cubilia=amet&curae=diam&mauris=in&viverra=magna&diam=bibendum&vitae=imperdiet&quam=nullam&suspendisse=orci&potenti=pede&nullam=venenatis&porttitor=non&lacus=sodales&at=sed&turpis=tincidunt&donec=eu&posuere=felis&metus=fusce&vitae=posuere&ipsum=felis&aliquam=sed&non=lacus&mauris=morbi&morbi=sem&non=mauris&lectus=laoreet&aliquam=ut'
$ws.Range("I12").Value = 'https://pcworld.com/in.json?suspendisse=arcu&accumsan=sed&tortor=augue&quis=aliquam&turpis=erat&sed=volutpat&ante=in&vivamus=congue&tortor=etiam&duis=justo&mattis=etiam&egestas=pretium&metus=iaculis&aenean=justo&fermentum=in&donec=hac&ut=habitasse&mauris=platea&eget=dictumst&massa=etiam&tempor=faucibus&convallis=cursus&nulla=urna&neque=ut&libero=tellus&convallis=nulla&eget=ut&eleifend=erat&luctus=id&ultricies=mauris&eu=vulputate&nibh=elementum&quisque=nullam&id=varius&justo=nulla&sit=facilisi&amet=cras&sapien=non&dignissim=velit&vestibulum=nec&vestibulum=nisi&ante=vulputate&ipsum=nonummy&primis=maecenas&in=tincidunt&faucibus=lacus&orci=at&luctus=velit&et=vivamus&ultrices=vel&posuere=nulla&cubilia=eget&curae=eros&nulla=elementum&dapibus=pellentesque&dolor=quisque&vel=porta&est=volutpat&donec=erat&odio=quisque&justo=erat&sollicitudin=eros&ut=viverra&suscipit=eget&a=congue&feugiat=eget&et=semper&eros=rutrum&vestibulum=nulla&ac=nunc&est=purus&lacinia=phasellus&nisi=in&venenatis=felis&tristique=donec'
$ws.Range("I13").Value = 'http://ted.com/ullamcorper/augue/a/suscipit.js?nulla=ipsum&ac=integer&enim=a&in=nibh&tempor=in&turpis=quis&nec=justo&euismod=maecenas&scelerisque=rhoncus&quam=aliquam&turpis=lacus&adipiscing=morbi&lorem=quis&vitae=tortor&mattis=id&nibh=nulla&ligula=ultrices&nec=aliquet&sem=maecenas&duis=leo&aliquam=odio&convallis=condimentum&nunc=id&proin=luctus&at=nec&turpis=molestie&a=sed&pede=justo&posuere=pellentesque&nonummy=viverra&integer=pede&non=ac&velit=diam'
$ws.Range("I14").Value = 'https://tumblr.com/consequat/nulla/nisl/nunc.json?ut=dictumst&dolor=morbi&morbi=vestibulum&vel=velit&lectus=id&in=pretium&quam=iaculis&fringilla=diam&rhoncus=erat&mauris=fermentum&enim=justo&leo=nec&rhoncus=condimentum&sed=neque&vestibulum=sapien&sit=placerat&amet=ante&cursus=nulla&id=justo&turpis=aliquam&integer=quis&aliquet=turpis&massa=eget&id=elit&lobortis=sodales&convallis=scelerisque&tortor=mauris&risus=sit&dapibus=amet&augue=eros&vel=suspendisse&accumsan=accumsan&tellus=tortor&nisi=quis&eu=turpis&orci=sed&mauris=ante'
$ws.Range("I15").Value = 'https://artisteer.com/accumsan/odio/curabitur/convallis/duis.jsp?cubilia=quis&curae=orci&donec=eget&pharetra=orci&magna=vehicula&vestibulum=condimentum&aliquet=curabitur&ultrices=in&erat=libero&tortor=ut&sollicitudin=massa&mi=volutpat&sit=convallis&amet=morbi&lobortis=odio&sapien=odio&sapien=elementum&non=eu&mi=interdum&integer=eu&ac=tincidunt&neque=in'
$ws.Range("I16").Value = 'http://yale.edu/turpis/adipiscing.json?odio=dictumst&curabitur=maecenas&convallis=ut&duis=massa&consequat=quis&dui=augue&nec=luctus&nisi=tincidunt&volutpat=nulla&eleifend=mollis&donec=molestie&ut=lorem&dolor=quisque&morbi=ut&vel=erat&lectus=curabitur&in=gravida&quam=nisi&fringilla=at&rhoncus=nibh&mauris=in&enim=hac&leo=habitasse&rhoncus=platea&sed=dictumst&vestibulum=aliquam&sit=augue&amet=quam&cursus=sollicitudin&id=vitae&turpis=consectetuer&integer=eget&aliquet=rutrum&massa=at&id=lorem&lobortis=integer&convallis=tincidunt&tortor=ante&risus=vel&dapibus=ipsum&augue=praesent&vel=blandit&accumsan=lacinia&tellus=erat&nisi=vestibulum&eu=sed&orci=magna&mauris=at&lacinia=nunc&sapien=commodo&quis=placerat&libero=praesent&nullam=blandit&sit=nam&amet=nulla&turpis=integer&elementum=pede&ligula=justo&vehicula=lacinia&consequat=eget&morbi=tincidunt&a=eget&ipsum=tempus&integer=vel&a=pede&nibh=morbi&in=porttitor'
$ws.Range("I17").Value = 'http://weibo.com/ultrices/phasellus/id/sapien/in/sapien/iaculis.json?habitasse=aliquet&platea=ultrices&dictumst=erat&aliquam=tortor&augue=sollicitudin&quam=mi&sollicitudin=sit&vitae=amet&consectetuer=lobortis&eget=sapien&rutrum=sapien&at=non&lorem=mi&integer=integer&tincidunt=ac&ante=neque&vel=duis&ipsum=bibendum&praesent=morbi&blandit=non&lacinia=quam&erat=nec&vestibulum=dui&sed=luctus&magna=rutrum&at=nulla&nunc=tellus&commodo=in&placerat=sagittis&praesent=dui&blandit=vel&nam=nisl&nulla=duis&integer=ac&pede=nibh&justo=fusce&lacinia=lacus&eget=purus&tincidunt=aliquet&eget=at&tempus=feugiat&vel=non&pede=pretium&morbi=quis&porttitor=lectus&lorem=suspendisse&id=potenti'
$ws.Range("I18").Value = 'https://newyorker.com/in/porttitor/pede/justo/eu/massa/donec.aspx?quisque=ipsum&erat=primis&eros=in&viverra=faucibus&eget=orci&congue=luctus&eget=et&semper=ultrices&rutrum=posuere&nulla=cubilia&nunc=curae&purus=mauris&phasellus=viverra&in=diam&felis=vitae&donec=quam&semper=suspendisse&sapien=potenti&a=nullam&libero=porttitor&nam=lacus&dui=at&proin=turpis&leo=donec&odio=posuere&porttitor=metus&id=vitae&consequat=ipsum&in=aliquam&consequat=non&ut=mauris&nulla=morbi&sed=non&accumsan=lectus&felis=aliquam&ut=sit&at=amet&dolor=diam'
$ws.Range("I19").Value = 'http://a8.net/libero/rutrum/ac/lobortis/vel/dapibus/at.jsp?in=orci&hac=eget&habitasse=orci&platea=vehicula&dictumst=condimentum&etiam=curabitur&faucibus=in&cursus=libero&urna=ut&ut=massa&tellus=volutpat&nulla=convallis&ut=morbi&erat=odio&id=odio&mauris=elementum&vulputate=eu&elementum=interdum&nullam=eu&varius=tincidunt&nulla=in&facilisi=leo&cras=maecenas&non=pulvinar&velit=lobortis&nec=est&nisi=phasellus&vulputate=sit&nonummy=amet&maecenas=erat&tincidunt=nulla&lacus=tempus&at=vivamus&velit=in&vivamus=felis&vel=eu&nulla=sapien&eget=cursus&eros=vestibulum&elementum=proin&pellentesque=eu&quisque=mi&porta=nulla&volutpat=ac&erat=enim&quisque=in&erat=tempor&eros=turpis&viverra=nec&eget=euismod&congue=scelerisque&eget=quam&semper=turpis&rutrum=adipiscing&nulla=lorem&nunc=vitae&purus=mattis&phasellus=nibh&in=ligula&felis=nec&donec=sem&semper=duis&sapien=aliquam&a=convallis&libero=nunc&nam=proin&dui=at&proin=turpis&leo=a&odio=pede&porttitor=posuere&id=nonummy&consequat=integer&in=non&consequat=velit&ut=donec&nulla=diam&sed=neque&accumsan=vestibulum&felis=eget'
$ws.Range("I20").Value = 'http://comcast.net/volutpat.json?semper=nulla&est=tellus&quam=in&pharetra=sagittis&magna=dui&ac=vel&consequat=nisl&metus=duis&sapien=ac&ut=nibh&nunc=fusce&vestibulum=lacus&ante=purus&ipsum=aliquet&primis=at&in=feugiat&faucibus=non&orci=pretium&luctus=quis&et=lectus&ultrices=suspendisse&posuere=potenti&cubilia=in&curae=eleifend&mauris=quam&viverra=a&diam=odio&vitae=in&quam=hac&suspendisse=habitasse&potenti=platea&nullam=dictumst&porttitor=maecenas&lacus=ut&at=massa&turpis=quis&donec=augue&posuere=luctus&metus=tincidunt&vitae=nulla&ipsum=mollis&aliquam=molestie&non=lorem&mauris=quisque&morbi=ut&non=erat&lectus=curabitur&aliquam=gravida&sit=nisi&amet=at&diam=nibh&in=in&magna=hac&bibendum=habitasse&imperdiet=platea&nullam=dictumst&orci=aliquam&pede=augue&venenatis=quam&non=sollicitudin&sodales=vitae&sed=consectetuer&tincidunt=eget&eu=rutrum&felis=at&fusce=lorem&posuere=integer&felis=tincidunt&sed=ante&lacus=vel&morbi=ipsum&sem=praesent&mauris=blandit&laoreet=lacinia&ut=erat&rhoncus=vestibulum&aliquet=sed&pulvinar=magna&sed=at&nisl=nunc&nunc=commodo&rhoncus=placerat&dui=praesent&vel=blandit&sem=nam&sed=nulla&sagittis=integer&nam=pede&congue=justo&risus=lacinia&semper=eget&porta=tincidunt&volutpat=eget'
$ws.Range("I21").Value = 'https://vk.com/blandit/nam/nulla/integer/pede.html?pede=ridiculus&morbi=mus&porttitor=vivamus&lorem=vestibulum&id=sagittis&ligula=sapien&suspendisse=cum&ornare=sociis&consequat=natoque&lectus=penatibus&in=et&est=magnis&risus=dis&auctor=parturient&sed=montes&tristique=nascetur&in=ridiculus&tempus=mus&sit=etiam&amet=vel&sem=augue&fusce=vestibulum&consequat=rutrum&nulla=rutrum&nisl=neque&nunc=aenean&nisl=auctor&duis=gravida&bibendum=sem&felis=praesent&sed=id&interdum=massa&venenatis=id&turpis=nisl&enim=venenatis&blandit=lacinia&mi=aenean&in=sit&porttitor=amet&pede=justo&justo=morbi&eu=ut&massa=odio&donec=cras&dapibus=mi&duis=pede&at=malesuada&velit=in&eu=imperdiet&est=et&congue=commodo&elementum=vulputate&in=justo'

$ws.Range("I21").Select()
